$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308 - this pushes existing rows 308-325 down
# to 309-326, preserving their data (matches target diff exactly).
$ws.Rows.Item(308).Insert()

# Populate the newly inserted row 308 with the new weekly record.
$ws.Cells.Item(308, 1).Value = 8
$ws.Cells.Item(308, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(308, 3).Value = "Coquimbo"
$ws.Cells.Item(308, 4).Value = 44706
$ws.Cells.Item(308, 5).Value = 4
$ws.Cells.Item(308, 6).Value = 100112032
$ws.Cells.Item(308, 7).Value = "Zapallo italiano"
$ws.Cells.Item(308, 8).Value = "Sin especificar"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 400
$ws.Cells.Item(308, 11).Value = 15000
$ws.Cells.Item(308, 12).Value = 15500
$ws.Cells.Item(308, 13).Value = 15250
$ws.Cells.Item(308, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(308, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(308, 16).Value = 254
$ws.Cells.Item(308, 17).Value = 60
$ws.Cells.Item(308, 18).Value = "Hortaliza"
